$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.853.79'
$ws.Range('E2').Value = '  +0.12%  '
$ws.Range('D3').Value = '3.520.99'
$ws.Range('E3').Value = '  +2.90%  '
$ws.Range('E4').Value = '  +0.13%  '
$ws.Range('D5').Value = "'591.80"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.37%  '
$ws.Range('D6').Value = "'136.43"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.09%  '
$ws.Range('D7').Value = '3.520.16'
$ws.Range('E7').Value = '  +2.86%  '
$ws.Range('E8').Value = '  +0.13%  '
$ws.Range('E9').Value = '  +1.34%  '
$ws.Range('E10').Value = '  +1.74%  '
$ws.Range('E11').Value = '  -1.89%  '
$ws.Range('D12').Value = "'0.382"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.79%  '
$ws.Range('D13').Value = '4.122.38'
$ws.Range('E13').Value = '  +3.32%  '
$ws.Range('B14').Value = 'ShibaInu'
$ws.Range('C14').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D14').Value = "'0.0000180"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.22%  '
$ws.Range('E15').Value = '  +2.37%  '
$ws.Range('B16').Value = 'WrappedEther'
$ws.Range('C16').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D16').Value = '3.528.66'
$ws.Range('E16').Value = '  +2.97%  '
$ws.Range('E17').Value = '  +1.51%  '
$ws.Range('D18').Value = '64.873.15'
$ws.Range('E18').Value = '  +0.40%  '
$ws.Range('D19').Value = "'9.95"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +4.02%  '
$ws.Range('E20').Value = '  -0.06%  '
$ws.Range('D21').Value = "'14.14"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +4.05%  '
$ws.Range('D22').Value = "'387.42"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.00%  '
$ws.Range('D23').Value = "'0.572"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +4.26%  '
$ws.Range('D24').Value = '3.665.70'
$ws.Range('E24').Value = '  +3.18%  '
$ws.Range('D25').Value = "'73.77"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +2.38%  '
$ws.Range('D26').Value = "'1.00"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.01%  '
$ws.Range('E27').Value = '  +6.69%  '
$ws.Range('D28').Value = "'7.62"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +5.71%  '
$ws.Range('E29').Value = '  -0.21%  '
$ws.Range('E30').Value = '  +3.04%  '
$ws.Range('D31').Value = "'8.16"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.33%  '
$ws.Range('D32').Value = '3.538.03'
$ws.Range('E32').Value = '  +3.35%  '
$ws.Range('E33').Value = '  +0.04%  '
$ws.Range('D34').Value = "'23.62"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +3.06%  '
$ws.Range('D35').Value = "'1.34"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +13.31%  '
$ws.Range('E36').Value = '  +1.26%  '
$ws.Range('D37').Value = "'170.21"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +1.31%  '
$ws.Range('E38').Value = '  +6.09%  '
$ws.Range('E39').Value = '  +1.07%  '
$ws.Range('E40').Value = '  +6.38%  '
$ws.Range('D41').Value = "'0.0795"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +5.01%  '
$ws.Range('E42').Value = '  +1.09%  '
$ws.Range('D43').Value = "'26.38"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +16.69%  '
$ws.Range('E44').Value = '  +0.11%  '
$ws.Range('D45').Value = "'42.42"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +1.50%  '
$ws.Range('E46').Value = '  +2.10%  '
$ws.Range('D47').Value = "'1.20"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +5.46%  '
$ws.Range('E48').Value = '  +1.98%  '
$ws.Range('D49').Value = "'6.83"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +5.86%  '
$ws.Range('D50').Value = '2.395.17'
$ws.Range('E50').Value = '  +10.62%  '
$ws.Range('D51').Value = "'302.88"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +10.68%  '
